# Update NATMI TPM LR-pair data (Mst1-Mst1r) for YoungD0:
# - delete the three "MuSCs" sending-cluster rows (old rows 8-10)
# - recompute/update the remaining rows 2-7 with the new TPM-derived values
# - dimension shrinks from A1:T10 to A1:T7 automatically as a result

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the trailing three data rows (old rows 8,9,10 -- "MuSCs" sender rows)
$ws.Range("A8:T10").EntireRow.Delete()

# Row 2 : ECs -> FAPs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Mst1"
$ws.Range("C2").Value = "Mst1r"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.280701
$ws.Range("H2").Value = 0.842103
$ws.Range("I2").Value = 0.188797685202
$ws.Range("J2").Value = 0.188797685202
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 5.171219666666667
$ws.Range("N2").Value = 15.513659
$ws.Range("O2").Value = 0.8031574816043436
$ws.Range("P2").Value = 0.8031574816043435
$ws.Range("Q2").Value = 1.451566531653
$ws.Range("R2").Value = 13.064098784877
$ws.Range("S2").Value = 0.151634273379568
$ws.Range("T2").Value = 0.151634273379568

# Row 3 : ECs -> MuSCs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Mst1"
$ws.Range("C3").Value = "Mst1r"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.280701
$ws.Range("H3").Value = 0.842103
$ws.Range("I3").Value = 0.188797685202
$ws.Range("J3").Value = 0.188797685202
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.267392666666667
$ws.Range("N3").Value = 3.802178
$ws.Range("O3").Value = 0.1968425183956564
$ws.Range("P3").Value = 0.1968425183956563
$ws.Range("Q3").Value = 0.355758388926
$ws.Range("R3").Value = 3.201825500334
$ws.Range("S3").Value = 0.03716341182243203
$ws.Range("T3").Value = 0.03716341182243203

# Row 4 : FAPs -> FAPs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Mst1"
$ws.Range("C4").Value = "Mst1r"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.9248883333333334
$ws.Range("H4").Value = 2.774665
$ws.Range("I4").Value = 0.6220739377617791
$ws.Range("J4").Value = 0.622073937761779
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 5.171219666666667
$ws.Range("N4").Value = 15.513659
$ws.Range("O4").Value = 0.8031574816043436
$ws.Range("P4").Value = 0.8031574816043435
$ws.Range("Q4").Value = 4.782800738803889
$ws.Range("R4").Value = 43.04520664923501
$ws.Range("S4").Value = 0.4996233372244477
$ws.Range("T4").Value = 0.4996233372244476

# Row 5 : FAPs -> MuSCs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Mst1"
$ws.Range("C5").Value = "Mst1r"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.9248883333333334
$ws.Range("H5").Value = 2.774665
$ws.Range("I5").Value = 0.6220739377617791
$ws.Range("J5").Value = 0.622073937761779
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.267392666666667
$ws.Range("N5").Value = 3.802178
$ws.Range("O5").Value = 0.1968425183956564
$ws.Range("P5").Value = 0.1968425183956563
$ws.Range("Q5").Value = 1.172196691152222
$ws.Range("R5").Value = 10.54977022037
$ws.Range("S5").Value = 0.1224506005373314
$ws.Range("T5").Value = 0.1224506005373314

# Row 6 : MuSCs -> FAPs
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Mst1"
$ws.Range("C6").Value = "Mst1r"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.2811926666666666
$ws.Range("H6").Value = 0.8435779999999999
$ws.Range("I6").Value = 0.1891283770362209
$ws.Range("J6").Value = 0.1891283770362209
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 5.171219666666667
$ws.Range("N6").Value = 15.513659
$ws.Range("O6").Value = 0.8031574816043436
$ws.Range("P6").Value = 0.8031574816043435
$ws.Range("Q6").Value = 1.454109047989111
$ws.Range("R6").Value = 13.086981431902
$ws.Range("S6").Value = 0.151899871000328
$ws.Range("T6").Value = 0.151899871000328

# Row 7 : MuSCs -> MuSCs
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Mst1"
$ws.Range("C7").Value = "Mst1r"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.2811926666666666
$ws.Range("H7").Value = 0.8435779999999999
$ws.Range("I7").Value = 0.1891283770362209
$ws.Range("J7").Value = 0.1891283770362209
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.267392666666667
$ws.Range("N7").Value = 3.802178
$ws.Range("O7").Value = 0.1968425183956564
$ws.Range("P7").Value = 0.1968425183956563
$ws.Range("Q7").Value = 0.3563815236537777
$ws.Range("R7").Value = 3.207433712883999
$ws.Range("S7").Value = 0.03722850603589295
$ws.Range("T7").Value = 0.03722850603589294
